$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row at row 10 ---
# This shifts current rows 10..20 down to 11..21.
$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").Value = 44827
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = 300000000
$ws.Range("G10").Value = "Espárragos"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 2200
$ws.Range("L10").Value = 2300
$ws.Range("M10").Value = 2250
$ws.Range("N10").Value = "`$/kilo"
$ws.Range("O10").Value = "Provincia de Diguillín"
$ws.Range("P10").Value = 2250
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = "Hortaliza"

# --- Insert second new row at row 16 ---
# This shifts current rows 16..21 down to 17..22.
$ws.Rows.Item(16).Insert()

$ws.Range("A16").Value = 7
$ws.Range("B16").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C16").Value = "Ñuble"
$ws.Range("D16").Value = 44460
$ws.Range("E16").Value = 16
$ws.Range("F16").Value = 300000000
$ws.Range("G16").Value = "Espárragos"
$ws.Range("H16").Value = "Verde"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = 2200
$ws.Range("L16").Value = 2300
$ws.Range("M16").Value = 2250
$ws.Range("N16").Value = "`$/kilo"
$ws.Range("O16").Value = "Provincia de Diguillín"
$ws.Range("P16").Value = 2250
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = "Hortaliza"

$wb.Save()
